$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-03 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-04 Thursday", 2)
$d.Content.Find.Execute("51÷2=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "96÷5=19, 1", 2)
$d.Content.Find.Execute("93÷3=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "91÷4=22, 3", 2)
$d.Content.Find.Execute("35÷4=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "61÷4=15, 1", 2)
$d.Content.Find.Execute("89÷9=9, 8", $true, $false, $false, $false, $false, $true, 1, $false, "74÷5=14, 4", 2)
$d.Content.Find.Execute("77÷2=38, 1", $true, $false, $false, $false, $false, $true, 1, $false, "27÷9=3, 0", 2)
$d.Content.Find.Execute("90÷2=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=13, 3", 2)
$d.Content.Find.Execute("58÷2=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "71÷6=11, 5", 2)
$d.Content.Find.Execute("70÷4=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "45÷2=22, 1", 2)
$d.Content.Find.Execute("25÷5=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2)
$d.Content.Find.Execute("91÷2=45, 1", $true, $false, $false, $false, $false, $true, 1, $false, "17÷3=5, 2", 2)
$d.Content.Find.Execute("43÷8=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "88÷2=44, 0", 2)
$d.Content.Find.Execute("23÷2=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "88÷8=11, 0", 2)
$d.Content.Find.Execute("69÷9=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=34, 1", 2)
$d.Content.Find.Execute("88÷9=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "82÷5=16, 2", 2)
$d.Content.Find.Execute("84÷8=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=22, 1", 2)
$d.Content.Find.Execute("14÷5=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "45÷4=11, 1", 2)
$d.Content.Find.Execute("28÷5=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "26÷4=6, 2", 2)
$d.Content.Find.Execute("14÷6=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "69÷5=13, 4", 2)
$d.Content.Find.Execute("56÷6=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "11÷4=2, 3", 2)
$d.Content.Find.Execute("52÷7=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "70÷3=23, 1", 2)
$d.Content.Find.Execute("10÷9=1, 1", $true, $false, $false, $false, $false, $true, 1, $false, "86÷2=43, 0", 2)
$d.Content.Find.Execute("80÷3=26, 2", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=29, 1", 2)
$d.Content.Find.Execute("27÷3=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "36÷4=9, 0", 2)
$d.Content.Find.Execute("22÷4=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "37÷6=6, 1", 2)
$d.Content.Find.Execute("35÷3=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "29÷9=3, 2", 2)
